$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row cells: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
$ws.Range("A1").Value2 = "Segmentname_FV2404"
$ws.Range("B1").Value2 = "Segmentgruppe_FV2404"
$ws.Range("C1").Value2 = "Segment_FV2404"
$ws.Range("D1").Value2 = "Datenelement_FV2404"
$ws.Range("E1").Value2 = "Segment ID_FV2404"
$ws.Range("F1").Value2 = "Code_FV2404"
$ws.Range("G1").Value2 = "Qualifier_FV2404"
$ws.Range("H1").Value2 = "Beschreibung_FV2404"
$ws.Range("I1").Value2 = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value2 = "Bedingung_FV2404"

$ws.Range("L1").Value2 = "Segmentname_FV2410"
$ws.Range("M1").Value2 = "Segmentgruppe_FV2410"
$ws.Range("N1").Value2 = "Segment_FV2410"
$ws.Range("O1").Value2 = "Datenelement_FV2410"
$ws.Range("P1").Value2 = "Segment ID_FV2410"
$ws.Range("Q1").Value2 = "Code_FV2410"
$ws.Range("R1").Value2 = "Qualifier_FV2410"
$ws.Range("S1").Value2 = "Beschreibung_FV2410"
$ws.Range("T1").Value2 = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value2 = "Bedingung_FV2410"

# Turn the data range into an Excel Table (ListObject)
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U65"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# Freeze the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
